$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.938819938647526
$ws.Range("B3").Value = 8.595745181530894
$ws.Range("B4").Value = 8.379316793265259
$ws.Range("B5").Value = 8.289805315195812
$ws.Range("B6").Value = 8.274866895163106
$ws.Range("B7").Value = 8.378114741901969
$ws.Range("B8").Value = 8.821805501660975
$ws.Range("B9").Value = 9.640695799583801
$ws.Range("B10").Value = 10.20495719013445
$ws.Range("B11").Value = 10.45248734535772
$ws.Range("B12").Value = 10.54483155020264
$ws.Range("B13").Value = 10.52500646403383
$ws.Range("B14").Value = 10.46011280418878
$ws.Range("B15").Value = 10.42018051987005
$ws.Range("B16").Value = 10.18858976415794
$ws.Range("B17").Value = 10.04411585700317
$ws.Range("B18").Value = 9.960160273863568
$ws.Range("B19").Value = 9.931589401009361
$ws.Range("B20").Value = 10.05958471727975
$ws.Range("B21").Value = 10.47921190877326
$ws.Range("B22").Value = 10.7453327207005
$ws.Range("B23").Value = 10.60406468796972
$ws.Range("B24").Value = 10.05259402946669
$ws.Range("B25").Value = 9.425332191843514
$ws.Range("C2").Value = 5.600988739745013
$ws.Range("C3").Value = 5.344437272379126
$ws.Range("C4").Value = 5.179321475052843
$ws.Range("C5").Value = 5.110174352179095
$ws.Range("C6").Value = 5.098581682427193
$ws.Range("C7").Value = 5.178396398896196
$ws.Range("C8").Value = 5.51413192611831
$ws.Range("C9").Value = 6.110539792136334
$ws.Range("C10").Value = 6.509086533959637
$ws.Range("C11").Value = 6.681506434958139
$ws.Range("C12").Value = 6.745502237611032
$ws.Range("C13").Value = 6.731777491317749
$ws.Range("C14").Value = 6.686797474023811
$ws.Range("C15").Value = 6.659076658829463
$ws.Range("C16").Value = 6.497638087373522
$ws.Range("C17").Value = 6.396311535119288
$ws.Range("C18").Value = 6.337196124392731
$ws.Range("C19").Value = 6.317037874434843
$ws.Range("C20").Value = 6.40718449283596
$ws.Range("C21").Value = 6.70004449668768
$ws.Range("C22").Value = 6.883885859364282
$ws.Range("C23").Value = 6.786463254071489
$ws.Range("C24").Value = 6.402271510404864
$ws.Range("C25").Value = 5.956043802802055
$ws.Range("D2").Value = 4.572609394117377
$ws.Range("D3").Value = 4.515740681239266
$ws.Range("D4").Value = 4.479922716352376
$ws.Range("D5").Value = 4.465107973894142
$ws.Range("D6").Value = 4.462635023427376
$ws.Range("D7").Value = 4.479723794386929
$ws.Range("D8").Value = 4.553192520359437
$ws.Range("D9").Value = 4.689791849368124
$ws.Range("D10").Value = 4.785174179689648
$ws.Range("D11").Value = 4.827401096023093
$ws.Range("D12").Value = 4.843217351457426
$ws.Range("D13").Value = 4.839818890874811
$ws.Range("D14").Value = 4.828705837955122
$ws.Range("D15").Value = 4.821875890675936
$ws.Range("D16").Value = 4.782390550488811
$ws.Range("D17").Value = 4.757864685058442
$ws.Range("D18").Value = 4.743648986238735
$ws.Range("D19").Value = 4.738817277864347
$ws.Range("D20").Value = 4.760486844689381
$ws.Range("D21").Value = 4.831974791406906
$ws.Range("D22").Value = 4.87767814767852
$ws.Range("D23").Value = 4.853380778849004
$ws.Range("D24").Value = 4.759301725385157
$ws.Range("D25").Value = 4.653682617797074
$ws.Range("E2").Value = 16.48800361265329
$ws.Range("E3").Value = 15.55441097959251
$ws.Range("E4").Value = 14.95632413581002
$ws.Range("E5").Value = 14.70661821184786
$ws.Range("E6").Value = 14.66480244666628
$ws.Range("E7").Value = 14.95298033449859
$ws.Range("E8").Value = 16.171406248977
$ws.Range("E9").Value = 18.43668239167298
$ws.Range("E10").Value = 20.07305316876849
$ws.Range("E11").Value = 20.7752704493664
$ws.Range("E12").Value = 21.03515283316829
$ws.Range("E13").Value = 20.979450269646
$ws.Range("E14").Value = 20.79677183394269
$ws.Range("E15").Value = 20.68409143438198
$ws.Range("E16").Value = 20.02631573808457
$ws.Range("E17").Value = 19.61201277992648
$ws.Range("E18").Value = 19.36974492644785
$ws.Range("E19").Value = 19.28703403892656
$ws.Range("E20").Value = 19.65652685998974
$ws.Range("E21").Value = 20.85059236686358
$ws.Range("E22").Value = 21.5958387922888
$ws.Range("E23").Value = 21.20129168996221
$ws.Range("E24").Value = 19.63641475436781
$ws.Range("E25").Value = 17.79645522314502
$ws.Range("F2").Value = 22.86306017478282
$ws.Range("F3").Value = 22.85798319484641
$ws.Range("F4").Value = 22.8615069846806
$ws.Range("F5").Value = 22.86461198251124
$ws.Range("F6").Value = 22.86522831809003
$ws.Range("F7").Value = 22.86154210416057
$ws.Range("F8").Value = 22.85993109341566
$ws.Range("F9").Value = 22.90945244221206
$ws.Range("F10").Value = 22.97785299817613
$ws.Range("F11").Value = 23.01587354065716
$ws.Range("F12").Value = 23.0312574314405
$ws.Range("F13").Value = 23.02790048672481
$ws.Range("F14").Value = 23.01711945740166
$ws.Range("F15").Value = 23.01064399768282
$ws.Range("F16").Value = 22.97550667843961
$ws.Range("F17").Value = 22.95571531780283
$ws.Range("F18").Value = 22.94498224192505
$ws.Range("F19").Value = 22.9414600897933
$ws.Range("F20").Value = 22.95775487203276
$ws.Range("F21").Value = 23.02025939722623
$ws.Range("F22").Value = 23.066854839142
$ws.Range("F23").Value = 23.04146281555333
$ws.Range("F24").Value = 22.95683077914503
$ws.Range("F25").Value = 22.89042059074052
$ws.Range("I2").Value = 19.31656217468171
$ws.Range("I3").Value = 19.39331494579942
$ws.Range("I4").Value = 19.44440815652837
$ws.Range("I5").Value = 19.46622464860703
$ws.Range("I6").Value = 19.46990733749204
$ws.Range("I7").Value = 19.44469835274089
$ws.Range("I8").Value = 19.3422021219401
$ws.Range("I9").Value = 19.17276468284186
$ws.Range("I10").Value = 19.06763133117378
$ws.Range("I11").Value = 19.0240302538781
$ws.Range("I12").Value = 19.00812924808032
$ws.Range("I13").Value = 19.01152665793708
$ws.Range("I14").Value = 19.02270983468863
$ws.Range("I15").Value = 19.02963933183471
$ws.Range("I16").Value = 19.07056596924025
$ws.Range("I17").Value = 19.0967567846798
$ws.Range("I18").Value = 19.11221855301081
$ws.Range("I19").Value = 19.11752184425536
$ws.Range("I20").Value = 19.09392757091886
$ws.Range("I21").Value = 19.01940849728153
$ws.Range("I22").Value = 18.97426127101496
$ws.Range("I23").Value = 18.99803113275983
$ws.Range("I24").Value = 19.09520539915011
$ws.Range("I25").Value = 19.21521179681097
$ws.Range("K2").Value = 8.736392222137528
$ws.Range("K3").Value = 8.491567826991453
$ws.Range("K4").Value = 8.339027531667458
$ws.Range("K5").Value = 8.276404596909028
$ws.Range("K6").Value = 8.265981147734605
$ws.Range("K7").Value = 8.338184714561923
$ws.Range("K8").Value = 8.652492726037943
$ws.Range("K9").Value = 9.247343204631596
$ws.Range("K10").Value = 9.666472813777041
$ws.Range("K11").Value = 9.852373386948615
$ws.Range("K12").Value = 9.922019169569593
$ws.Range("K13").Value = 9.907054073074653
$ws.Range("K14").Value = 9.858118587352095
$ws.Range("K15").Value = 9.828044578586272
$ws.Range("K16").Value = 9.654222155501122
$ws.Range("K17").Value = 9.546318880176299
$ws.Range("K18").Value = 9.483810919519184
$ws.Range("K19").Value = 9.462572501044857
$ws.Range("K20").Value = 9.557851904534138
$ws.Range("K21").Value = 9.872512993061866
$ws.Range("K22").Value = 10.11412956890199
$ws.Range("K23").Value = 9.975358314987377
$ws.Range("K24").Value = 9.552639292211651
$ws.Range("K25").Value = 9.089243263582297
$ws.Range("N2").Value = 17.52924621832931
$ws.Range("N3").Value = 17.58768085849675
$ws.Range("N4").Value = 17.62520077763389
$ws.Range("N5").Value = 17.64090436391063
$ws.Range("N6").Value = 17.64353697064661
$ws.Range("N7").Value = 17.62541088403113
$ws.Range("N8").Value = 17.54905472874963
$ws.Range("N9").Value = 17.4122784791467
$ws.Range("N10").Value = 17.3196013855352
$ws.Range("N11").Value = 17.27911795372297
$ws.Range("N12").Value = 17.26402753495288
$ws.Range("N13").Value = 17.26726688153117
$ws.Range("N14").Value = 17.27787165724953
$ws.Range("N15").Value = 17.28439857916186
$ws.Range("N16").Value = 17.32228069194867
$ws.Range("N17").Value = 17.34594850363385
$ws.Range("N18").Value = 17.35971941260347
$ws.Range("N19").Value = 17.36440914665342
$ws.Range("N20").Value = 17.34341270071323
$ws.Range("N21").Value = 17.2747502796066
$ws.Range("N22").Value = 17.23127247654346
$ws.Range("N23").Value = 17.25434996333096
$ws.Range("N24").Value = 17.34455862633682
$ws.Range("N25").Value = 17.44790188080151
$ws.Range("O2").Value = 20.43177991603502
$ws.Range("O3").Value = 20.48199027864781
$ws.Range("O4").Value = 20.51773012443045
$ws.Range("O5").Value = 20.53352538254923
$ws.Range("O6").Value = 20.53622240527001
$ws.Range("O7").Value = 20.51793816586435
$ws.Range("O8").Value = 20.44807106865914
$ws.Range("O9").Value = 20.35018428811926
$ws.Range("O10").Value = 20.30232048843285
$ws.Range("O11").Value = 20.28580655244911
$ws.Range("O12").Value = 20.28031199204047
$ws.Range("O13").Value = 20.28146155369951
$ws.Range("O14").Value = 20.28533928839144
$ws.Range("O15").Value = 20.28781341615472
$ws.Range("O16").Value = 20.30350577803943
$ws.Range("O17").Value = 20.31448158417405
$ws.Range("O18").Value = 20.32128949743259
$ws.Range("O19").Value = 20.32367945823811
$ws.Range("O20").Value = 20.31326194860177
$ws.Range("O21").Value = 20.28417968822002
$ws.Range("O22").Value = 20.26959720436036
$ws.Range("O23").Value = 20.27697454701338
$ws.Range("O24").Value = 20.31381179543953
$ws.Range("O25").Value = 20.37245473735766
